$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the old "code commit" note from F3
$ws.Range("F3").ClearContents()

# Remove old "update" and "ici is getting update..." notes from F7/G8
$ws.Range("F7").ClearContents()
$ws.Range("G8").ClearContents()

# Add the new note replacing them, now placed at F5
$ws.Range("F5").Value = "validate collectionsName in server"

# Update the active selection to match the new state
$ws.Range("M10").Select()
